# Adds three new LeetCode tracking rows (270, 272, 2411) to Sheet1 and
# updates the sheet view (zoom/selection) to match the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 107: Smallest Subarrays With Maximum Bitwise OR -------------------
$ws.Range("A107").Value = 2411
$ws.Range("B107").Value = "Smallest Subarrays With Maximum Bitwise OR"
$ws.Range("C107").Value = "#array #bit-minipulation #sliding-window "
$ws.Range("D107").Value = "medium"
$ws.Range("E107").Value = 0
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 50
$ws.Range("H107").Value = 45867
$ws.Range("H107").NumberFormat = "m/d/yy"
$ws.Range("I107").Value = 45867
$ws.Range("I107").NumberFormat = "m/d/yy"

# --- Row 108: Closest Binary Search Tree Value II ---------------------------
$ws.Range("A108").Value = 272
$ws.Range("B108").Value = "Closest Binary Search Tree Value II"
$ws.Range("C108").Value = "#two-pointers #stack #tree #dfs #bst #heap #binary-tree "
$ws.Range("D108").Value = "hard"
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = 1
$ws.Range("G108").Value = 15
$ws.Range("H108").Value = 45867
$ws.Range("H108").NumberFormat = "m/d/yy"
$ws.Range("I108").Value = 45867
$ws.Range("I108").NumberFormat = "m/d/yy"
$ws.Range("J108").Value = "deque"

# --- Row 109: Closest Binary Search Tree Value ------------------------------
$ws.Range("A109").Value = 270
$ws.Range("B109").Value = "Closest Binary Search Tree Value"
$ws.Range("C109").Value = "#tree #bst "
$ws.Range("D109").Value = "easy"
$ws.Range("E109").Value = 1
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 15
$ws.Range("H109").Value = 45867
$ws.Range("H109").NumberFormat = "m/d/yy"
$ws.Range("I109").Value = 45867
$ws.Range("I109").NumberFormat = "m/d/yy"

# --- Row heights to match the source formatting -----------------------------
$ws.Rows.Item(107).RowHeight = 51
$ws.Rows.Item(108).RowHeight = 68
$ws.Rows.Item(109).RowHeight = 34

# --- Sheet view: zoom + scroll position + active selection ------------------
$ws.Range("H109:I109").Select()
$excel.ActiveWindow.Zoom = 172
